# Change the year in the astromap link from 2018 to 2022.
#
# The paragraph looks like:
#   (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/).
# split across three differently-formatted runs: "(", the hyperlink-styled
# URL, and ").". The edit collapses that into plain text with the updated
# year and no special character formatting (no Hyperlink style, no custom
# fonts/size).

$d = $word.ActiveDocument

$newText = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

foreach ($p in $d.Paragraphs) {
    $pRange = $p.Range
    if ($pRange.Text -like "*amper.ped.muni.cz*GaNight/2018*") {
        # Range over the paragraph's text, excluding the trailing paragraph mark.
        $rng = $d.Range($pRange.Start, $pRange.End - 1)

        # Replace the whole (formatted, multi-run) link text with plain,
        # unformatted text in one go.
        $rng.Delete()
        $rng.InsertAfter($newText)
    }
}
